# Fill in the "Actual time length to complete" column for the two tasks
# that were just finished (DQ2 response 2 / row 14, DQ2 response 3 / row 15).
# Values are stored as Excel time-serial fractions of a day:
#   C14 = 30 minutes = 0.5 / 24
#   C15 = 15 minutes = 0.25 / 24
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C14").Value = 0.020833333333333332
$ws.Range("C15").Value = 0.010416666666666666

# The Total row (C20 = SUM(C2:C19)) recalculates automatically.

# Move the active selection down to the next task to work on.
$ws.Range("C16").Select()
